$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.691559791564941
$ws.Range("B1").Value = 2.835700511932373
$ws.Range("C1").Value = 1.407265305519104
$ws.Range("D1").Value = 0.6675335764884949
$ws.Range("E1").Value = 0.6232852935791016
